$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "26÷8=3, 2"
$t.Cell(1, 2).Range.Text = "30÷3=10, 0"
$t.Cell(1, 3).Range.Text = "17÷6=2, 5"
$t.Cell(1, 4).Range.Text = "75÷5=15, 0"
$t.Cell(1, 5).Range.Text = "27÷6=4, 3"

$t.Cell(5, 1).Range.Text = "89÷9=9, 8"
$t.Cell(5, 2).Range.Text = "15÷3=5, 0"
$t.Cell(5, 3).Range.Text = "50÷4=12, 2"
$t.Cell(5, 4).Range.Text = "72÷8=9, 0"
$t.Cell(5, 5).Range.Text = "36÷6=6, 0"

$t.Cell(9, 1).Range.Text = "36÷9=4, 0"
$t.Cell(9, 2).Range.Text = "62÷3=20, 2"
$t.Cell(9, 3).Range.Text = "65÷8=8, 1"
$t.Cell(9, 4).Range.Text = "54÷8=6, 6"
$t.Cell(9, 5).Range.Text = "21÷8=2, 5"

$t.Cell(13, 1).Range.Text = "76÷5=15, 1"
$t.Cell(13, 2).Range.Text = "83÷4=20, 3"
$t.Cell(13, 3).Range.Text = "63÷6=10, 3"
$t.Cell(13, 4).Range.Text = "16÷8=2, 0"
$t.Cell(13, 5).Range.Text = "69÷5=13, 4"

$t.Cell(17, 1).Range.Text = "66÷5=13, 1"
$t.Cell(17, 2).Range.Text = "50÷2=25, 0"
$t.Cell(17, 3).Range.Text = "61÷7=8, 5"
$t.Cell(17, 4).Range.Text = "36÷3=12, 0"
$t.Cell(17, 5).Range.Text = "19÷6=3, 1"

